# DOMA-1155 regenerate excel templates
# Adds "Unit" column after "Address" and splits the single "Reading" column
# into four tariff-specific reading columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at C for "Unit" (shifts Service..Source right by 1)
$ws.Columns("C:C").Insert()

# 2) Insert three more columns at H:J for the extra tariff readings
#    (after the original "Reading" column, now at G, shifts Contact/Source right)
$ws.Columns("H:J").Insert()

# Row 1 - headers
$ws.Range("C1").Value = "Unit"
$ws.Range("G1").Value = "Reading from tariff №1"
$ws.Range("H1").Value = "Reading from tariff №2"
$ws.Range("I1").Value = "Reading from tariff №3"
$ws.Range("J1").Value = "Reading from tariff №4"

# Row 2 - template placeholders
$ws.Range("C2").Value = "{d.meter[i].unitName}"
$ws.Range("H2").Value = "{d.meter[i].value2}"
$ws.Range("I2").Value = "{d.meter[i].value3}"
$ws.Range("J2").Value = "{d.meter[i].value4}"

# Row 3 - template placeholders
$ws.Range("C3").Value = "{d.meter[i + 1].unitName}"
$ws.Range("H3").Value = "{d.meter[i + 1].value2}"
$ws.Range("I3").Value = "{d.meter[i + 1].value3}"
$ws.Range("J3").Value = "{d.meter[i + 1].value4}"
